$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ref, [string]$val) {
    $cell = $ws.Range($ref)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '27.043.23'
$ws.Range('E2').Value = '  +0.43%  '

# Row 3
$ws.Range('D3').Value = '1.677.59'
$ws.Range('E3').Value = '  +0.47%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
Set-CellText 'D5' '215.26'
$ws.Range('E5').Value = '  -0.03%  '

# Row 6
$ws.Range('E6').Value = '  -0.21%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('E8').Value = '  +1.94%  '

# Row 9
Set-CellText 'D9' '21.33'
$ws.Range('E9').Value = '  +5.20%  '

# Row 10
$ws.Range('E10').Value = '  +0.14%  '

# Row 11
Set-CellText 'D11' '0.0885'
$ws.Range('E11').Value = '  -0.34%  '

# Row 12
$ws.Range('D12').Value = '1.914.16'
$ws.Range('E12').Value = '  +0.44%  '

# Row 13
$ws.Range('D13').Value = '1.681.60'
$ws.Range('E13').Value = '  +0.65%  '

# Row 14
Set-CellText 'D14' '4.13'
$ws.Range('E14').Value = '  +1.19%  '

# Row 15
Set-CellText 'D15' '0.537'
$ws.Range('E15').Value = '  +1.89%  '

# Row 16
Set-CellText 'D16' '66.39'
$ws.Range('E16').Value = '  +1.26%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '27.030.79'
$ws.Range('E17').Value = '  +0.39%  '

# Row 18
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText 'D18' '8.20'
$ws.Range('E18').Value = '  +3.08%  '

# Row 19
Set-CellText 'D19' '237.77'
$ws.Range('E19').Value = '  +1.58%  '

# Row 20
$ws.Range('E20').Value = '  +0.41%  '

# Row 21
$ws.Range('E21').Value = '  +0.08%  '

# Row 22
$ws.Range('E22').Value = '  +1.88%  '

# Row 23
Set-CellText 'D23' '9.31'
$ws.Range('E23').Value = '  +1.81%  '

# Row 24
Set-CellText 'D24' '2.12'
$ws.Range('E24').Value = '  -1.30%  '

# Row 25
Set-CellText 'D25' '146.51'
$ws.Range('E25').Value = '  +0.19%  '

# Row 26
Set-CellText 'D26' '7.23'
$ws.Range('E26').Value = '  +1.81%  '

# Row 27
Set-CellText 'D27' '16.40'
$ws.Range('E27').Value = '  +3.22%  '

# Row 28
$ws.Range('E28').Value = '  +1.08%  '

# Row 29
Set-CellText 'D29' '0.999'
$ws.Range('E29').Value = '  -0.18%  '

# Row 30
$ws.Range('E30').Value = '  +0.41%  '

# Row 31
$ws.Range('E31').Value = '  +0.07%  '

# Row 32
$ws.Range('E32').Value = '  +0.81%  '

# Row 33
$ws.Range('D33').Value = '1.539.17'
$ws.Range('E33').Value = '  +5.94%  '

# Row 34
$ws.Range('E34').Value = '  +1.63%  '

# Row 35
$ws.Range('E35').Value = '  +3.27%  '

# Row 36
$ws.Range('E36').Value = '  -1.28%  '

# Row 37
Set-CellText 'D37' '0.595'
$ws.Range('E37').Value = '  +1.69%  '

# Row 38
Set-CellText 'D38' '0.916'
$ws.Range('E38').Value = '  +1.64%  '

# Row 39
$ws.Range('E39').Value = '  +2.21%  '

# Row 40
$ws.Range('E40').Value = '  +3.14%  '

# Row 41
$ws.Range('E41').Value = '  +0.05%  '

# Row 42
Set-CellText 'D42' '67.62'
$ws.Range('E42').Value = '  +1.97%  '

# Row 43
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText 'D43' '2.27'
$ws.Range('E43').Value = '  -1.33%  '

# Row 44
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText 'D44' '5.51'
$ws.Range('E44').Value = '  -4.03%  '

# Row 45
$ws.Range('D45').Value = '1.821.11'
$ws.Range('E45').Value = '  +0.62%  '

# Row 46
$ws.Range('E46').Value = '  +0.22%  '

# Row 47
Set-CellText 'D47' '90.74'
$ws.Range('E47').Value = '  +0.28%  '

# Row 48
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText 'D48' '1.55'
$ws.Range('E48').Value = '  +0.92%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText 'D49' '0.104'
$ws.Range('E49').Value = '  +2.12%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D50' '8.04'
$ws.Range('E50').Value = '  +5.72%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText 'D51' '0.0509'
$ws.Range('E51').Value = '  +0.48%  '
